$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.751.95"
$ws.Range("E2").Value = "  +1.31%  "

$ws.Range("D3").Value = "1.658.48"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("D4").Formula = "=""0.9999"""
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Formula = "=""1.000"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("D6").Formula = "=""303.90"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").Formula = "=""0.3814"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  +0.70%  "

$ws.Range("D8").Formula = "=""0.3631"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Formula = "=""51.18"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -1.36%  "

$ws.Range("D10").Formula = "=""1.249"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +1.57%  "

$ws.Range("D11").Formula = "=""0.08219"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").Formula = "=""1.000"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Formula = "=""22.71"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +1.08%  "

$ws.Range("D14").Formula = "=""6.537"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +0.81%  "

$ws.Range("D15").Formula = "=""7.460"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +1.09%  "

$ws.Range("D16").Formula = "=""0.00001236"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").Value = "1.660.78"
$ws.Range("E17").Value = "  +1.68%  "

$ws.Range("D18").Formula = "=""97.61"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +2.47%  "

$ws.Range("D19").Formula = "=""0.07014"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +0.99%  "

$ws.Range("D20").Formula = "=""6.821"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +3.38%  "

$ws.Range("D21").Formula = "=""17.77"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +1.39%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Formula = "=""12.90"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +2.78%  "

$ws.Range("D24").Value = "23.744.37"
$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("D25").Formula = "=""2.526"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +0.27%  "

$ws.Range("D26").Formula = "=""3.062"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -0.52%  "

$ws.Range("D27").Formula = "=""21.33"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +0.63%  "

$ws.Range("D28").Formula = "=""152.66"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").Formula = "=""5.214"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("D30").Formula = "=""134.85"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +0.73%  "

$ws.Range("D31").Value = "1.844.21"
$ws.Range("E31").Value = "  +1.61%  "

$ws.Range("D32").Formula = "=""7.072"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +6.52%  "

$ws.Range("D33").Formula = "=""2.214"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +2.75%  "

$ws.Range("D34").Formula = "=""1.072"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +1.87%  "

$ws.Range("D35").Formula = "=""11.74"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +3.46%  "

$ws.Range("D36").Formula = "=""0.02826"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +2.42%  "

$ws.Range("D37").Formula = "=""0.2538"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +1.63%  "

$ws.Range("D38").Formula = "=""6.136"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +1.46%  "

$ws.Range("D39").Formula = "=""0.08803"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +0.21%  "

$ws.Range("D40").Formula = "=""0.07117"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("D41").Formula = "=""13.09"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +7.03%  "

$ws.Range("D42").Formula = "=""0.7065"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("D43").Formula = "=""1.336"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -0.28%  "

$ws.Range("D44").Formula = "=""16.17"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)

$ws.Range("D45").Formula = "=""0.6540"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").Formula = "=""2.326"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +2.08%  "

$ws.Range("E47").Value = "  +0.14%  "

$ws.Range("D48").Formula = "=""3.981"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +0.30%  "

$ws.Range("D49").Formula = "=""0.07959"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -0.65%  "

$ws.Range("D50").Formula = "=""128.81"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("D51").Formula = "=""1.191"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -0.34%  "

$excel.CutCopyMode = $false
